$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the targetCode for row 4 from "NPC_roras" to the new "MOB_Slime" string
$ws.Range("C4").Value = "MOB_Slime"

# Reset selection to C4 (was F4)
$ws.Range("C4").Select()

# Clear the distinguishing cell borders/styles on A3, A4, A5, B2 so they fall back
# to the plain default style (matches the style-table cleanup in the diff).
$ws.Range("B2").Borders.LineStyle = 0
$ws.Range("A3").Borders.LineStyle = 0
$ws.Range("A4").Borders.LineStyle = 0
$ws.Range("A5").Borders.LineStyle = 0
